# Implements unmodified logit for reliability based capacity expansion:
# switch the "ETLE-capacity" sheet's logit exponent (B2) from -8 to -0.3,
# and make that sheet the active/selected tab (it was left open/selected
# when the author saved the workbook).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ETLE-capacity")
$ws.Activate()
$ws.Range("B2").Value = -0.3
